# Generate Report for Handoff
# The localization CI run regenerated the handoff package under a new
# GUID-named markdown file (ab856b28-4594-45c5-9549-62b71948049b.md,
# replacing 738aeb2f-6057-4f34-bf1c-de92e9382f94.md) and re-ran the
# handoff/xliff generation a little later. Update the status workbook
# cells (and the matching hyperlink display text) to reflect this,
# without touching the still-valid external hyperlink target URLs.

$wb = $excel.ActiveWorkbook

$oldGuid = "738aeb2f-6057-4f34-bf1c-de92e9382f94"
$newGuid = "ab856b28-4594-45c5-9549-62b71948049b"
$oldHash = "626913f026d6e992421b722ffb58205ef69923bb"
$newHash = "a7337b20b74832f7024b12f7462e7350dddaad6a"
$oldHyperlinkAddress = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a7ea708ef7d3473d7b59cba4772d5d671be28049/e2e/$oldGuid.md"

# ---- Overview sheet ----
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = "$newGuid.md"
$wsOverview.Range("B2").Value = "e2e\$newGuid.md"
$wsOverview.Range("G2").Value = "2016-08-26 16:58:15"

# Keep the existing external link target, just refresh the displayed text.
$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), $oldHyperlinkAddress, [Type]::Missing, [Type]::Missing, "e2e\$newGuid.md")

# ---- zh-cn sheet ----
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("A2").Value = "$newGuid.md"
$wsZhCn.Range("G2").Value = "$newGuid.$newHash.zh-cn.xlf"
$wsZhCn.Range("H2").Value = "2016-08-26 16:58:09"

$wsZhCn.Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), $oldHyperlinkAddress, [Type]::Missing, [Type]::Missing, "$newGuid.md")

# ---- de-de sheet ----
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("A2").Value = "$newGuid.md"
$wsDeDe.Range("G2").Value = "$newGuid.$newHash.de-de.xlf"
$wsDeDe.Range("H2").Value = "2016-08-26 16:58:15"

$wsDeDe.Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), $oldHyperlinkAddress, [Type]::Missing, [Type]::Missing, "$newGuid.md")
